{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the async (context) => { ... } function.\n//\n// Goal: after the last paragraph (\"para los mensajes \"), append:\n//   1. an empty paragraph\n//   2. \"php artisan make:model Secretaria -mcr \"\n//   3. \"para crear de una vez el modelo, controlador y las migraciones del modulo secretaria\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document always ends with the \"para los mensajes \" paragraph just\n// before the sectPr, so anchor on the last paragraph in the body.\nconst anchorParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// 1) Blank paragraph right after \"para los mensajes \".\nconst blankParagraph = anchorParagraph.insertParagraph(\"\", \"After\");\n\n// 2) Command line paragraph.\nconst commandParagraph = blankParagraph.insertParagraph(\n  \"php artisan make:model Secretaria -mcr \",\n  \"After\"\n);\n\n// 3) Explanation paragraph.\ncommandParagraph.insertParagraph(\n  \"para crear de una vez el modelo, controlador y las migraciones del modulo secretaria\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Goal: after the last paragraph (\"para los mensajes \"), append:\n#   1. an empty paragraph\n#   2. \"php artisan make:model Secretaria -mcr \"\n#   3. \"para crear de una vez el modelo, controlador y las migraciones del modulo secretaria\"\n\n$d = $word.ActiveDocument\n\n# The document always ends with the \"para los mensajes \" paragraph just\n# before the sectPr, so anchor on the last paragraph in the body.\n$anchor = $d.Paragraphs.Last\n$anchor.Range.InsertParagraphAfter()\n\n# 1) Blank paragraph is now the last paragraph; add the command line text\n#    as a new paragraph right after it.\n$blank = $d.Paragraphs.Last\n$blank.Range.InsertParagraphAfter()\n\n$commandPara = $d.Paragraphs.Last\n$commandPara.Range.InsertAfter(\"php artisan make:model Secretaria -mcr \")\n$commandPara.Range.InsertParagraphAfter()\n\n# 3) Explanation paragraph.\n$explainPara = $d.Paragraphs.Last\n$explainPara.Range.InsertAfter(\"para crear de una vez el modelo, controlador y las migraciones del modulo secretaria\")\n"}
